$d = $word.ActiveDocument

# Locate the "Impact" bullet list that lives under the
# "KEY ACHIEVEMENTS AND IMPACT" heading. We anchor on paragraph position
# (heading index + 2, since the very next paragraph is the "Impact"
# sub-heading) rather than matching bullet text directly, because some of
# the bullet strings in this list are duplicated verbatim elsewhere in the
# resume (e.g. under "Partner - Siege Analytics") and a text-based
# Find/Replace could touch the wrong occurrence.
$headingIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $headingIndex = $i
    }
}

if ($headingIndex -eq -1) {
    throw "Could not locate 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

$firstBulletIndex = $headingIndex + 2   # heading, "Impact" sub-heading, then bullets

# Capture references to the six existing bullet paragraphs up front so that
# later deletions don't shift the indices of paragraphs we still need to
# edit.
$p1 = $d.Paragraphs($firstBulletIndex)       # "Discovered systematic race coding errors..."
$p2 = $d.Paragraphs($firstBulletIndex + 1)   # "Achieved 87% prediction accuracy..."
$p3 = $d.Paragraphs($firstBulletIndex + 2)   # "Built redistricting platform..."
$p4 = $d.Paragraphs($firstBulletIndex + 3)   # "Developed longitudinal data analysis methods..." (removed)
$p5 = $d.Paragraphs($firstBulletIndex + 4)   # "Provided expert testimony..." (removed)
$p6 = $d.Paragraphs($firstBulletIndex + 5)   # "Demystified FEC compliance..."

$p1.Range.Text = "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
$p2.Range.Text = "• Reduced polling margins from ±4.2% to ±2.1%"
$p3.Range.Text = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
$p6.Range.Text = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# Remove the two bullets that have no replacement (delete higher index
# first so the still-pending deletion's Range stays valid).
$p5.Range.Delete()
$p4.Range.Delete()
